$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto data rows 2-51 (columns B=Coin, C=Link, D=Price, E=Volume(1h))
$data = @(
    ,@(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '30.547.54', '  +1.22%  ')
    ,@(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.850.58', '  +0.08%  ')
    ,@(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.000', '  -0.07%  ')
    ,@(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '233.38', '  +0.33%  ')
    ,@(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.000', '  -0.01%  ')
    ,@(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4719', '  +0.99%  ')
    ,@(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2739', '  +1.24%  ')
    ,@(9, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.06311', '  -0.78%  ')
    ,@(10, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '17.68', '  +9.24%  ')
    ,@(11, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.848.95', '  -0.15%  ')
    ,@(12, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07453', '  +0.54%  ')
    ,@(13, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.003', '  +1.41%  ')
    ,@(14, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '84.35', '  -0.57%  ')
    ,@(15, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.6244', '  -0.30%  ')
    ,@(16, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '30.495.42', '  +1.16%  ')
    ,@(17, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '243.39', '  +7.41%  ')
    ,@(18, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.000', '  +0.01%  ')
    ,@(19, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '12.64', '  +0.70%  ')
    ,@(20, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000007328', '  +0.25%  ')
    ,@(21, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.000', '  -0.11%  ')
    ,@(22, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '4.917', '  -0.28%  ')
    ,@(23, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '5.924', '  +0.13%  ')
    ,@(24, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '9.158', '  -0.70%  ')
    ,@(25, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '162.69', '  -2.01%  ')
    ,@(26, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '17.94', '  +1.23%  ')
    ,@(27, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.875', '  +0.79%  ')
    ,@(28, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1020', '  -1.49%  ')
    ,@(29, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.361', '  -1.56%  ')
    ,@(30, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.010', '  -2.19%  ')
    ,@(31, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '3.820', '  -0.87%  ')
    ,@(32, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.04837', '  -0.81%  ')
    ,@(33, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.134', '  -1.67%  ')
    ,@(34, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7002', '  -1.61%  ')
    ,@(35, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.711', '  +0.65%  ')
    ,@(36, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01896', '  +2.58%  ')
    ,@(37, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.679', '  +1.63%  ')
    ,@(38, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.8724', '  -3.56%  ')
    ,@(39, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '1.984', '  +1.76%  ')
    ,@(40, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '106.35', '  +1.22%  ')
    ,@(41, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.001', '  +0.24%  ')
    ,@(42, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '5.519', '  -0.28%  ')
    ,@(43, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.4054', '  -0.28%  ')
    ,@(44, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '7.168', '  +2.14%  ')
    ,@(45, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '62.37', '  +3.99%  ')
    ,@(46, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1211', '  +2.20%  ')
    ,@(47, 'Elrond', 'https://coinranking.com/coin/omwkOTglq+elrond-egld', '33.43', '  +1.23%  ')
    ,@(48, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '8.543', '  -0.22%  ')
    ,@(49, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.05540', '  -0.54%  ')
    ,@(50, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '1.353', '  -2.09%  ')
    ,@(51, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.3656', '  -0.09%  ')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = $row[4]
}
